# New Submission Synced: 2026-02-11 06:28:42
# Sheet "JSS 3E" gets a new form response appended as row 13, and the
# existing row 12 "Admission No" (C12), which had been stored as text,
# is corrected to a real number — matching how every other row in the
# column is typed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")
$ws.Activate()

# --- fix C12: was inline text "44" -> numeric 44 -------------------------
$ws.Range("C12").Value = 44

# --- append row 13: the newly synced submission ---------------------------
$ws.Range("A13").Value = "2026-02-11 06:28:42"
$ws.Range("B13").Value = "HAUWA ALIYU MADUGU"

# C13 ("Admission No") must remain a text value "37" (same quirk as the
# original sheet, where some admission numbers are stored as text).
# A plain Value/Formula assignment of "37" gets auto-coerced to the number
# 37, and forcing it via NumberFormat="@" leaves a stray style behind, so
# instead build it as a text formula result and flatten it to a literal
# value via copy / paste-special (values only).
$ws.Range("C13").Formula = '=TEXT(37,"0")'
$ws.Range("C13").Copy()
$ws.Range("C13").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("D13").Value = 10
